# Add a new worksheet "Sheet7" at the end of the workbook (after the last
# existing sheet) and make it the active sheet/tab.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add($null, $lastSheet)
$ws7.Name = "Sheet7"

# --- Content: "With hidden rows" example table -----------------------------
# (values are entered in the same order the original author typed them, so
# that new shared-string entries line up with the reference workbook)
$ws7.Range("A1").Value = "With hidden rows"

$ws7.Range("A3").Value = "Table 1"

$ws7.Range("A4").Value = "a1"
$ws7.Range("A7").Value = "a4"

$ws7.Range("B4").Value = "b1"
$ws7.Range("B7").Value = "b4"

$ws7.Range("C4").Value = "c1"
$ws7.Range("C7").Value = "c4"

$ws7.Range("A5").Value = "a2"
$ws7.Range("B5").Value = "b2"
$ws7.Range("C5").Value = "c2"

$ws7.Range("A6").Value = "a3"
$ws7.Range("B6").Value = "b3"
$ws7.Range("C6").Value = "c3"

# --- Content: color examples ------------------------------------------------
$ws7.Range("E1").Value = "And colors:"

$ws7.Range("H1").Value = "(200,201,202)"
$ws7.Range("I1").Value = "pattern"
$ws7.Range("J1").Value = "no color"
$ws7.Range("G1").Value = "{'theme':5}"

# Apply the fills in the same order the original workbook recorded them
# (theme color, then plain RGB color, then a hatched pattern fill).
$ws7.Range("G1").Interior.ThemeColor = 6
$ws7.Range("H1").Interior.Color = 13289928
$ws7.Range("I1").Interior.Pattern = 12

# --- Hidden rows / columns ---------------------------------------------------
$ws7.Range("B1").ColumnWidth = -0.75
$ws7.Range("B1").EntireColumn.Hidden = $true

$ws7.Range("A5:A6").EntireRow.Hidden = $true

$ws7.Range("G1").ColumnWidth = 9.59
$ws7.Range("H1").ColumnWidth = 11.65

# --- Page setup / view -------------------------------------------------------
$ws7.PageSetup.Orientation = 1

$ws7.Range("C7").Select()
